$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for four new log entries (rows 48-51) right after the
#    existing last entry (row 47). Inserting copies down formatting
#    (number formats, alignment, shared-formula continuation, and it
#    also keeps the old "Total Hours" summary row intact, shifting it
#    down and auto-adjusting its SUM() range).
# ------------------------------------------------------------------
$ws.Rows("48:51").Insert()

# The summary row (previously row 50) is now row 54. The new entries
# should end up at rows 48-51 and the summary row should move further
# down to row 57 (5 blank rows between the last entry and the total),
# matching the final layout of the log.
$ws.Rows("54:56").Insert()

# ------------------------------------------------------------------
# 2. New entry: 7 Jul 2022, 08:00-08:30, Code, Segformer paper review
# ------------------------------------------------------------------
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 44749
$ws.Range("C48").Value = 0.33333333333333331
$ws.Range("D48").Value = 0.35416666666666669
$ws.Range("E48").Formula = "=D48-C48"
$ws.Range("F48").Value = "Code"
$ws.Range("G48").Value = "1. Segformer paper review -  half done"

# ------------------------------------------------------------------
# 3. New entry: 8 Jul 2022, 08:15-09:00, Code, mix-FFN/DS conv/Prenorm
# ------------------------------------------------------------------
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 44750
$ws.Range("C49").Value = 0.34375
$ws.Range("D49").Value = 0.375
$ws.Range("E49").Formula = "=D49-C49"
$ws.Range("F49").Value = "Code"
$ws.Range("G49").Value = "1. Mix feedforward, DS conv, Prenorm layers understanding"

# ------------------------------------------------------------------
# 4. New entry: 9 Jul 2022, 13:30-15:30, Code, self-attention/einops
# ------------------------------------------------------------------
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 44751
$ws.Range("C50").Value = 0.5625
$ws.Range("D50").Value = 0.64583333333333337
$ws.Range("E50").Formula = "=D50-C50"
$ws.Range("F50").Value = "Code"
$ws.Range("G50").Value = "1. Efficient self attention understanding" + [char]10 + "2. einops rearrange understanding"
$ws.Rows(50).RowHeight = 30

# ------------------------------------------------------------------
# 5. New entry: 11 Jul 2022, 08:00-08:30, Code, MiT module
# ------------------------------------------------------------------
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 44753
$ws.Range("C51").Value = 0.33333333333333331
$ws.Range("D51").Value = 0.35416666666666669
$ws.Range("E51").Formula = "=D51-C51"
$ws.Range("F51").Value = "Code"
$ws.Range("G51").Value = "1. MiT module working understanding - half"

# ------------------------------------------------------------------
# 6. Update the view so the newly added rows are visible / selected,
#    matching the author's final cursor position.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A52").Select() | Out-Null

$excel.Calculate() | Out-Null
